# Auto-generated cryptos.xlsx update script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'30.412.13"
$ws.Range('E2').Value = '  -0.10%  '

# Row 3
$ws.Range('D3').Value = "'1.927.35"
$ws.Range('E3').Value = '  +4.07%  '

# Row 4
$ws.Range('D4').Value = "'0.9989"
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').Value = "'240.08"
$ws.Range('E5').Value = '  +2.89%  '

# Row 6
$ws.Range('D6').Value = "'0.9990"

# Row 7
$ws.Range('D7').Value = "'0.4767"
$ws.Range('E7').Value = '  +0.45%  '

# Row 8
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = "'0.2871"
$ws.Range('E8').Value = '  +4.42%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = "'0.06578"
$ws.Range('E9').Value = '  +3.83%  '

# Row 10
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = "'19.08"
$ws.Range('E10').Value = '  +7.83%  '

# Row 11
$ws.Range('B11').Value = 'Litecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D11').Value = "'107.52"
$ws.Range('E11').Value = '  +26.86%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = "'1.923.04"
$ws.Range('E12').Value = '  +2.07%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = "'0.07625"
$ws.Range('E13').Value = '  +2.37%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'5.146"
$ws.Range('E14').Value = '  +3.80%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = "'0.6605"
$ws.Range('E15').Value = '  +5.65%  '

# Row 16
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = "'308.28"
$ws.Range('E16').Value = '  +25.12%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = "'30.429.28"
$ws.Range('E17').Value = '  +0.10%  '

# Row 18
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = "'13.02"
$ws.Range('E18').Value = '  +2.71%  '

# Row 19
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = "'0.9990"
$ws.Range('E19').Value = '  -0.11%  '

# Row 20
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = "'0.000007534"
$ws.Range('E20').Value = '  +2.61%  '

# Row 21
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = "'2.166.40"
$ws.Range('E21').Value = '  +3.16%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'5.312"
$ws.Range('E22').Value = '  +8.34%  '

# Row 23
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = "'0.9992"
$ws.Range('E23').Value = '  -0.03%  '

# Row 24
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = "'6.289"
$ws.Range('E24').Value = '  +6.53%  '

# Row 25
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'167.67"
$ws.Range('E25').Value = '  +1.53%  '

# Row 26
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = "'9.235"
$ws.Range('E26').Value = '  +1.52%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'20.20"
$ws.Range('E27').Value = '  +12.32%  '

# Row 28
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = "'2.038"
$ws.Range('E28').Value = '  +8.61%  '

# Row 29
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = "'0.1117"
$ws.Range('E29').Value = '  +8.59%  '

# Row 30
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = "'1.358"
$ws.Range('E30').Value = '  +0.78%  '

# Row 31
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = "'4.102"
$ws.Range('E31').Value = '  +1.43%  '

# Row 32
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'3.930"
$ws.Range('E32').Value = '  +2.79%  '

# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.05014"
$ws.Range('E33').Value = '  +3.42%  '

# Row 34
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = "'0.7431"
$ws.Range('E34').Value = '  +6.39%  '

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'1.148"
$ws.Range('E35').Value = '  +1.58%  '

# Row 36
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.753"
$ws.Range('E36').Value = '  +1.54%  '

# Row 37
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = "'0.9982"
$ws.Range('E37').Value = '  -0.12%  '

# Row 38
$ws.Range('E38').Value = '  +2.68%  '

# Row 39
$ws.Range('D39').Value = "'2.705"
$ws.Range('E39').Value = '  +0.81%  '

# Row 40
$ws.Range('D40').Value = "'2.051"
$ws.Range('E40').Value = '  +2.78%  '

# Row 41
$ws.Range('D41').Value = "'0.8788"
$ws.Range('E41').Value = '  +0.18%  '

# Row 42
$ws.Range('D42').Value = "'106.89"
$ws.Range('E42').Value = '  +0.08%  '

# Row 43
$ws.Range('D43').Value = "'70.71"
$ws.Range('E43').Value = '  +11.53%  '

# Row 44
$ws.Range('D44').Value = "'5.806"
$ws.Range('E44').Value = '  +5.45%  '

# Row 45
$ws.Range('D45').Value = "'0.9985"
$ws.Range('E45').Value = '  -0.16%  '

# Row 46
$ws.Range('D46').Value = "'0.4166"
$ws.Range('E46').Value = '  +2.37%  '

# Row 47
$ws.Range('D47').Value = "'7.280"
$ws.Range('E47').Value = '  +1.47%  '

# Row 48
$ws.Range('D48').Value = "'9.267"
$ws.Range('E48').Value = '  +8.60%  '

# Row 49
$ws.Range('D49').Value = "'34.91"
$ws.Range('E49').Value = '  +2.95%  '

# Row 50
$ws.Range('D50').Value = "'0.1204"
$ws.Range('E50').Value = '  +0.56%  '

# Row 51
$ws.Range('D51').Value = "'0.05620"
$ws.Range('E51').Value = '  +2.00%  '
